# Apply the "Modificações para o paper" edit:
#  - B5 input value changes from 20 to 60 (recalculates K5 automatically)
#  - Active selection on Sheet1 ends up at K5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 60

$ws.Range("K5").Select()
